# Update odds for rows 2-6 (values refreshed by the data feed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19

$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.3
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 6
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("AA3").Value = 17
$ws.Range("AH3").Value = 11
$ws.Range("AI3").Value = 26
$ws.Range("AO3").Value = 9.5
$ws.Range("AT3").Value = 2.38
$ws.Range("AX3").Value = 34

$ws.Range("I4").Value = 2.7
$ws.Range("L4").Value = 3.6
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 13
$ws.Range("AE4").Value = 21
$ws.Range("AL4").Value = 29
$ws.Range("AW4").Value = 4.33
$ws.Range("AX4").Value = 17
$ws.Range("AY4").Value = 34

$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53

$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("U6").Value = 2.5
$ws.Range("V6").Value = 1.5
$ws.Range("AC6").Value = 7
$ws.Range("AQ6").Value = 29
$ws.Range("BA6").Value = 251

# Row 7 (lzKILqFr, Columbus Crew vs New York Red Bulls, 19:45) is gone from the
# feed; deleting it shifts the old row 8 (CYtcF2g2, Real Salt Lake vs
# Minnesota United) up into row 7, and the used range shrinks to A1:BD7.
$ws.Rows.Item(7).Delete()

# Refresh the odds on the row that is now row 7 (formerly row 8) to the
# latest feed values.
$ws.Range("G7").Value = 2.05
$ws.Range("I7").Value = 3.6
$ws.Range("J7").Value = 2.6
$ws.Range("L7").Value = 3.75
$ws.Range("O7").Value = 1.2
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.67
$ws.Range("R7").Value = 2.15
$ws.Range("W7").Value = 10
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 19
$ws.Range("AF7").Value = 34
$ws.Range("AH7").Value = 13
$ws.Range("AJ7").Value = 12
$ws.Range("AO7").Value = 11
$ws.Range("AR7").Value = 51
$ws.Range("AU7").Value = 7
$ws.Range("AW7").Value = 5.5
$ws.Range("AX7").Value = 17
